$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 6576985
$ws.Cells.Item(2, 5).Value = "Cercle Brugge"
$ws.Cells.Item(2, 6).Value = "Westerlo"
$ws.Cells.Item(2, 7).Value = 2
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = "H"
$ws.Cells.Item(2, 10).Value = 1.727
$ws.Cells.Item(2, 11).Value = 4.2
$ws.Cells.Item(2, 12).Value = 4.2
$ws.Cells.Item(2, 13).Value = 1.45
$ws.Cells.Item(2, 14).Value = 5
$ws.Cells.Item(2, 15).Value = 5.75
$ws.Cells.Item(2, 16).Value = -1.25
$ws.Cells.Item(2, 17).Value = 1.95
$ws.Cells.Item(2, 18).Value = 1.9
$ws.Cells.Item(2, 19).Value = 3.5
$ws.Cells.Item(2, 20).Value = 1.9
$ws.Cells.Item(2, 21).Value = 1.95
$ws.Cells.Item(2, 22).Value = 0.45
$ws.Cells.Item(2, 23).Value = -1
$ws.Cells.Item(2, 24).Value = -1
$ws.Cells.Item(2, 25).Value = 0.95
$ws.Cells.Item(2, 26).Value = -1
$ws.Cells.Item(2, 27).Value = -1
$ws.Cells.Item(2, 28).Value = 0.95

$ws.Cells.Item(3, 2).Value = 6576986
$ws.Cells.Item(3, 5).Value = "Gent"
$ws.Cells.Item(3, 6).Value = "Standard Liege"
$ws.Cells.Item(3, 7).Value = 3
$ws.Cells.Item(3, 8).Value = 1
$ws.Cells.Item(3, 9).Value = "H"
$ws.Cells.Item(3, 10).Value = 1.5
$ws.Cells.Item(3, 11).Value = 4.75
$ws.Cells.Item(3, 12).Value = 5.5
$ws.Cells.Item(3, 13).Value = 1.363
$ws.Cells.Item(3, 14).Value = 5.75
$ws.Cells.Item(3, 15).Value = 6.5
$ws.Cells.Item(3, 16).Value = -1.5
$ws.Cells.Item(3, 17).Value = 1.925
$ws.Cells.Item(3, 18).Value = 1.925
$ws.Cells.Item(3, 19).Value = 3.5
$ws.Cells.Item(3, 20).Value = 1.925
$ws.Cells.Item(3, 21).Value = 1.925
$ws.Cells.Item(3, 22).Value = 0.363
$ws.Cells.Item(3, 23).Value = -1
$ws.Cells.Item(3, 24).Value = -1
$ws.Cells.Item(3, 25).Value = 0.925
$ws.Cells.Item(3, 26).Value = -1
$ws.Cells.Item(3, 27).Value = 0.925
$ws.Cells.Item(3, 28).Value = -1

$ws.Cells.Item(26, 2).Value = 7030334
$ws.Cells.Item(26, 5).Value = "Cercle Brugge"
$ws.Cells.Item(26, 6).Value = "Genk"
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = 1
$ws.Cells.Item(26, 9).Value = "A"
$ws.Cells.Item(26, 10).Value = 2.75
$ws.Cells.Item(26, 11).Value = 3.5
$ws.Cells.Item(26, 12).Value = 2.25
$ws.Cells.Item(26, 13).Value = 2.4
$ws.Cells.Item(26, 14).Value = 3.5
$ws.Cells.Item(26, 15).Value = 2.55
$ws.Cells.Item(26, 16).Value = 0
$ws.Cells.Item(26, 17).Value = 1.85
$ws.Cells.Item(26, 18).Value = 2
$ws.Cells.Item(26, 19).Value = 3
$ws.Cells.Item(26, 20).Value = 1.9
$ws.Cells.Item(26, 21).Value = 1.95
$ws.Cells.Item(26, 22).Value = -1
$ws.Cells.Item(26, 23).Value = -1
$ws.Cells.Item(26, 24).Value = 1.55
$ws.Cells.Item(26, 25).Value = -1
$ws.Cells.Item(26, 26).Value = 1
$ws.Cells.Item(26, 27).Value = -1
$ws.Cells.Item(26, 28).Value = 0.95

$ws.Cells.Item(27, 2).Value = 6810007
$ws.Cells.Item(27, 5).Value = "Eupen"
$ws.Cells.Item(27, 6).Value = "Club Brugge"
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 5
$ws.Cells.Item(27, 9).Value = "A"
$ws.Cells.Item(27, 10).Value = 4.75
$ws.Cells.Item(27, 11).Value = 4
$ws.Cells.Item(27, 12).Value = 1.571
$ws.Cells.Item(27, 13).Value = 7
$ws.Cells.Item(27, 14).Value = 4.75
$ws.Cells.Item(27, 15).Value = 1.333
$ws.Cells.Item(27, 16).Value = 1.5
$ws.Cells.Item(27, 17).Value = 1.875
$ws.Cells.Item(27, 18).Value = 1.975
$ws.Cells.Item(27, 19).Value = 3
$ws.Cells.Item(27, 20).Value = 1.925
$ws.Cells.Item(27, 21).Value = 1.925
$ws.Cells.Item(27, 22).Value = -1
$ws.Cells.Item(27, 23).Value = -1
$ws.Cells.Item(27, 24).Value = 0.333
$ws.Cells.Item(27, 25).Value = -1
$ws.Cells.Item(27, 26).Value = 0.9750000000000001
$ws.Cells.Item(27, 27).Value = 0.925
$ws.Cells.Item(27, 28).Value = -1

$ws.Cells.Item(155, 2).Value = 6810130
$ws.Cells.Item(155, 5).Value = "Antwerp"
$ws.Cells.Item(155, 6).Value = "Westerlo"
$ws.Cells.Item(155, 7).Value = 2
$ws.Cells.Item(155, 8).Value = 2
$ws.Cells.Item(155, 9).Value = "D"
$ws.Cells.Item(155, 10).Value = 1.363
$ws.Cells.Item(155, 11).Value = 5
$ws.Cells.Item(155, 12).Value = 7.5
$ws.Cells.Item(155, 13).Value = 1.333
$ws.Cells.Item(155, 14).Value = 5.25
$ws.Cells.Item(155, 15).Value = 8
$ws.Cells.Item(155, 16).Value = -1.5
$ws.Cells.Item(155, 17).Value = 1.9
$ws.Cells.Item(155, 18).Value = 1.95
$ws.Cells.Item(155, 19).Value = 3
$ws.Cells.Item(155, 20).Value = 1.8
$ws.Cells.Item(155, 21).Value = 2.05
$ws.Cells.Item(155, 22).Value = -1
$ws.Cells.Item(155, 23).Value = 4.25
$ws.Cells.Item(155, 24).Value = -1
$ws.Cells.Item(155, 25).Value = -1
$ws.Cells.Item(155, 26).Value = 0.95
$ws.Cells.Item(155, 27).Value = 0.8
$ws.Cells.Item(155, 28).Value = -1

$ws.Cells.Item(156, 2).Value = 6810132
$ws.Cells.Item(156, 5).Value = "SintTruidense"
$ws.Cells.Item(156, 6).Value = "Charleroi"
$ws.Cells.Item(156, 7).Value = 1
$ws.Cells.Item(156, 8).Value = 0
$ws.Cells.Item(156, 9).Value = "H"
$ws.Cells.Item(156, 10).Value = 2.2
$ws.Cells.Item(156, 11).Value = 3.4
$ws.Cells.Item(156, 12).Value = 3.2
$ws.Cells.Item(156, 13).Value = 2.3
$ws.Cells.Item(156, 14).Value = 3.3
$ws.Cells.Item(156, 15).Value = 3
$ws.Cells.Item(156, 16).Value = -0.25
$ws.Cells.Item(156, 17).Value = 2
$ws.Cells.Item(156, 18).Value = 1.85
$ws.Cells.Item(156, 19).Value = 2.25
$ws.Cells.Item(156, 20).Value = 1.85
$ws.Cells.Item(156, 21).Value = 2
$ws.Cells.Item(156, 22).Value = 1.3
$ws.Cells.Item(156, 23).Value = -1
$ws.Cells.Item(156, 24).Value = -1
$ws.Cells.Item(156, 25).Value = 1
$ws.Cells.Item(156, 26).Value = -1
$ws.Cells.Item(156, 27).Value = -1
$ws.Cells.Item(156, 28).Value = 1

$ws.Cells.Item(164, 2).Value = 6810143
$ws.Cells.Item(164, 5).Value = "Charleroi"
$ws.Cells.Item(164, 6).Value = "KV Mechelen"
$ws.Cells.Item(164, 7).Value = 3
$ws.Cells.Item(164, 8).Value = 1
$ws.Cells.Item(164, 9).Value = "H"
$ws.Cells.Item(164, 10).Value = 2.05
$ws.Cells.Item(164, 11).Value = 3.4
$ws.Cells.Item(164, 12).Value = 3.5
$ws.Cells.Item(164, 13).Value = 2.375
$ws.Cells.Item(164, 14).Value = 3.1
$ws.Cells.Item(164, 15).Value = 3.1
$ws.Cells.Item(164, 16).Value = -0.25
$ws.Cells.Item(164, 17).Value = 2.025
$ws.Cells.Item(164, 18).Value = 1.825
$ws.Cells.Item(164, 19).Value = 2.25
$ws.Cells.Item(164, 20).Value = 2.025
$ws.Cells.Item(164, 21).Value = 1.825
$ws.Cells.Item(164, 22).Value = 1.375
$ws.Cells.Item(164, 23).Value = -1
$ws.Cells.Item(164, 24).Value = -1
$ws.Cells.Item(164, 25).Value = 1.025
$ws.Cells.Item(164, 26).Value = -1
$ws.Cells.Item(164, 27).Value = 1.025
$ws.Cells.Item(164, 28).Value = -1

$ws.Cells.Item(165, 2).Value = 6810139
$ws.Cells.Item(165, 5).Value = "Anderlecht"
$ws.Cells.Item(165, 6).Value = "Cercle Brugge"
$ws.Cells.Item(165, 7).Value = 2
$ws.Cells.Item(165, 8).Value = 0
$ws.Cells.Item(165, 9).Value = "H"
$ws.Cells.Item(165, 10).Value = 1.869
$ws.Cells.Item(165, 11).Value = 3.7
$ws.Cells.Item(165, 12).Value = 3.8
$ws.Cells.Item(165, 13).Value = 2.05
$ws.Cells.Item(165, 14).Value = 3.6
$ws.Cells.Item(165, 15).Value = 3.3
$ws.Cells.Item(165, 16).Value = -0.25
$ws.Cells.Item(165, 17).Value = 1.825
$ws.Cells.Item(165, 18).Value = 2.025
$ws.Cells.Item(165, 19).Value = 2.75
$ws.Cells.Item(165, 20).Value = 1.9
$ws.Cells.Item(165, 21).Value = 1.95
$ws.Cells.Item(165, 22).Value = 1.05
$ws.Cells.Item(165, 23).Value = -1
$ws.Cells.Item(165, 24).Value = -1
$ws.Cells.Item(165, 25).Value = 0.825
$ws.Cells.Item(165, 26).Value = -1
$ws.Cells.Item(165, 27).Value = -1
$ws.Cells.Item(165, 28).Value = 0.95

$ws.Cells.Item(181, 2).Value = 6810169
$ws.Cells.Item(181, 5).Value = "Westerlo"
$ws.Cells.Item(181, 6).Value = "Cercle Brugge"
$ws.Cells.Item(181, 7).Value = 4
$ws.Cells.Item(181, 8).Value = 2
$ws.Cells.Item(181, 9).Value = "H"
$ws.Cells.Item(181, 10).Value = 3.4
$ws.Cells.Item(181, 11).Value = 3.8
$ws.Cells.Item(181, 12).Value = 1.95
$ws.Cells.Item(181, 13).Value = 3.3
$ws.Cells.Item(181, 14).Value = 3.6
$ws.Cells.Item(181, 15).Value = 2.05
$ws.Cells.Item(181, 16).Value = 0.25
$ws.Cells.Item(181, 17).Value = 2
$ws.Cells.Item(181, 18).Value = 1.85
$ws.Cells.Item(181, 19).Value = 2.75
$ws.Cells.Item(181, 20).Value = 1.975
$ws.Cells.Item(181, 21).Value = 1.875
$ws.Cells.Item(181, 22).Value = 2.3
$ws.Cells.Item(181, 23).Value = -1
$ws.Cells.Item(181, 24).Value = -1
$ws.Cells.Item(181, 25).Value = 1
$ws.Cells.Item(181, 26).Value = -1
$ws.Cells.Item(181, 27).Value = 0.9750000000000001
$ws.Cells.Item(181, 28).Value = -1

$ws.Cells.Item(182, 2).Value = 6810167
$ws.Cells.Item(182, 5).Value = "Club Brugge"
$ws.Cells.Item(182, 6).Value = "KV Kortrijk"
$ws.Cells.Item(182, 7).Value = 3
$ws.Cells.Item(182, 8).Value = 3
$ws.Cells.Item(182, 9).Value = "D"
$ws.Cells.Item(182, 10).Value = 1.125
$ws.Cells.Item(182, 11).Value = 8.5
$ws.Cells.Item(182, 12).Value = 17
$ws.Cells.Item(182, 13).Value = 1.125
$ws.Cells.Item(182, 14).Value = 8.5
$ws.Cells.Item(182, 15).Value = 17
$ws.Cells.Item(182, 16).Value = -2.25
$ws.Cells.Item(182, 17).Value = 1.85
$ws.Cells.Item(182, 18).Value = 2
$ws.Cells.Item(182, 19).Value = 3.5
$ws.Cells.Item(182, 20).Value = 1.9
$ws.Cells.Item(182, 21).Value = 1.95
$ws.Cells.Item(182, 22).Value = -1
$ws.Cells.Item(182, 23).Value = 7.5
$ws.Cells.Item(182, 24).Value = -1
$ws.Cells.Item(182, 25).Value = -1
$ws.Cells.Item(182, 26).Value = 1
$ws.Cells.Item(182, 27).Value = 0.8999999999999999
$ws.Cells.Item(182, 28).Value = -1

$ws.Cells.Item(183, 2).Value = 6810168
$ws.Cells.Item(183, 5).Value = "OH Leuven"
$ws.Cells.Item(183, 6).Value = "Genk"
$ws.Cells.Item(183, 7).Value = 2
$ws.Cells.Item(183, 8).Value = 1
$ws.Cells.Item(183, 9).Value = "H"
$ws.Cells.Item(183, 10).Value = 4.5
$ws.Cells.Item(183, 11).Value = 4.2
$ws.Cells.Item(183, 12).Value = 1.666
$ws.Cells.Item(183, 13).Value = 4.333
$ws.Cells.Item(183, 14).Value = 4
$ws.Cells.Item(183, 15).Value = 1.7
$ws.Cells.Item(183, 16).Value = 0.75
$ws.Cells.Item(183, 17).Value = 1.95
$ws.Cells.Item(183, 18).Value = 1.9
$ws.Cells.Item(183, 19).Value = 3
$ws.Cells.Item(183, 20).Value = 1.975
$ws.Cells.Item(183, 21).Value = 1.875
$ws.Cells.Item(183, 22).Value = 3.333
$ws.Cells.Item(183, 23).Value = -1
$ws.Cells.Item(183, 24).Value = -1
$ws.Cells.Item(183, 25).Value = 0.95
$ws.Cells.Item(183, 26).Value = -1
$ws.Cells.Item(183, 27).Value = 0
$ws.Cells.Item(183, 28).Value = 0

$ws.Cells.Item(184, 2).Value = 6810165
$ws.Cells.Item(184, 5).Value = "Charleroi"
$ws.Cells.Item(184, 6).Value = "Eupen"
$ws.Cells.Item(184, 7).Value = 1
$ws.Cells.Item(184, 8).Value = 0
$ws.Cells.Item(184, 9).Value = "H"
$ws.Cells.Item(184, 10).Value = 1.6
$ws.Cells.Item(184, 11).Value = 4
$ws.Cells.Item(184, 12).Value = 5
$ws.Cells.Item(184, 13).Value = 1.8
$ws.Cells.Item(184, 14).Value = 3.75
$ws.Cells.Item(184, 15).Value = 4.2
$ws.Cells.Item(184, 16).Value = -0.75
$ws.Cells.Item(184, 17).Value = 2.05
$ws.Cells.Item(184, 18).Value = 1.8
$ws.Cells.Item(184, 19).Value = 2.75
$ws.Cells.Item(184, 20).Value = 1.95
$ws.Cells.Item(184, 21).Value = 1.9
$ws.Cells.Item(184, 22).Value = 0.8
$ws.Cells.Item(184, 23).Value = -1
$ws.Cells.Item(184, 24).Value = -1
$ws.Cells.Item(184, 25).Value = 0.5249999999999999
$ws.Cells.Item(184, 26).Value = -0.5
$ws.Cells.Item(184, 27).Value = -1
$ws.Cells.Item(184, 28).Value = 0.8999999999999999

$ws.Cells.Item(190, 2).Value = 6810171
$ws.Cells.Item(190, 5).Value = "KV Kortrijk"
$ws.Cells.Item(190, 6).Value = "Charleroi"
$ws.Cells.Item(190, 7).Value = 1
$ws.Cells.Item(190, 8).Value = 0
$ws.Cells.Item(190, 9).Value = "H"
$ws.Cells.Item(190, 10).Value = 3.2
$ws.Cells.Item(190, 11).Value = 3.5
$ws.Cells.Item(190, 12).Value = 2.1
$ws.Cells.Item(190, 13).Value = 3.4
$ws.Cells.Item(190, 14).Value = 3.4
$ws.Cells.Item(190, 15).Value = 2.05
$ws.Cells.Item(190, 16).Value = 0.25
$ws.Cells.Item(190, 17).Value = 2
$ws.Cells.Item(190, 18).Value = 1.85
$ws.Cells.Item(190, 19).Value = 2.5
$ws.Cells.Item(190, 20).Value = 1.925
$ws.Cells.Item(190, 21).Value = 1.925
$ws.Cells.Item(190, 22).Value = 2.4
$ws.Cells.Item(190, 23).Value = -1
$ws.Cells.Item(190, 24).Value = -1
$ws.Cells.Item(190, 25).Value = 1
$ws.Cells.Item(190, 26).Value = -1
$ws.Cells.Item(190, 27).Value = -1
$ws.Cells.Item(190, 28).Value = 0.925

$ws.Cells.Item(191, 2).Value = 6810174
$ws.Cells.Item(191, 5).Value = "Westerlo"
$ws.Cells.Item(191, 6).Value = "OH Leuven"
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 3
$ws.Cells.Item(191, 9).Value = "A"
$ws.Cells.Item(191, 10).Value = 1.909
$ws.Cells.Item(191, 11).Value = 3.75
$ws.Cells.Item(191, 12).Value = 3.5
$ws.Cells.Item(191, 13).Value = 1.909
$ws.Cells.Item(191, 14).Value = 3.5
$ws.Cells.Item(191, 15).Value = 3.8
$ws.Cells.Item(191, 16).Value = -0.5
$ws.Cells.Item(191, 17).Value = 1.925
$ws.Cells.Item(191, 18).Value = 1.925
$ws.Cells.Item(191, 19).Value = 2.5
$ws.Cells.Item(191, 20).Value = 1.85
$ws.Cells.Item(191, 21).Value = 2
$ws.Cells.Item(191, 22).Value = -1
$ws.Cells.Item(191, 23).Value = -1
$ws.Cells.Item(191, 24).Value = 2.8
$ws.Cells.Item(191, 25).Value = -1
$ws.Cells.Item(191, 26).Value = 0.925
$ws.Cells.Item(191, 27).Value = 0.8500000000000001
$ws.Cells.Item(191, 28).Value = -1

$ws.Cells.Item(241, 2).Value = 6942395
$ws.Cells.Item(241, 5).Value = "Gent"
$ws.Cells.Item(241, 6).Value = "Charleroi"
$ws.Cells.Item(241, 7).Value = 5
$ws.Cells.Item(241, 8).Value = 0
$ws.Cells.Item(241, 9).Value = "H"
$ws.Cells.Item(241, 10).Value = 1.571
$ws.Cells.Item(241, 11).Value = 4
$ws.Cells.Item(241, 12).Value = 5.75
$ws.Cells.Item(241, 13).Value = 1.4
$ws.Cells.Item(241, 14).Value = 4.333
$ws.Cells.Item(241, 15).Value = 8
$ws.Cells.Item(241, 16).Value = -1.25
$ws.Cells.Item(241, 17).Value = 2
$ws.Cells.Item(241, 18).Value = 1.85
$ws.Cells.Item(241, 19).Value = 2.75
$ws.Cells.Item(241, 20).Value = 1.95
$ws.Cells.Item(241, 21).Value = 1.9
$ws.Cells.Item(241, 22).Value = 0.3999999999999999
$ws.Cells.Item(241, 23).Value = -1
$ws.Cells.Item(241, 24).Value = -1
$ws.Cells.Item(241, 25).Value = 1
$ws.Cells.Item(241, 26).Value = -1
$ws.Cells.Item(241, 27).Value = 0.95
$ws.Cells.Item(241, 28).Value = -1

$ws.Cells.Item(242, 2).Value = 6810219
$ws.Cells.Item(242, 5).Value = "OH Leuven"
$ws.Cells.Item(242, 6).Value = "KV Mechelen"
$ws.Cells.Item(242, 7).Value = 1
$ws.Cells.Item(242, 8).Value = 0
$ws.Cells.Item(242, 9).Value = "H"
$ws.Cells.Item(242, 10).Value = 2.8
$ws.Cells.Item(242, 11).Value = 3.5
$ws.Cells.Item(242, 12).Value = 2.375
$ws.Cells.Item(242, 13).Value = 2.7
$ws.Cells.Item(242, 14).Value = 3.5
$ws.Cells.Item(242, 15).Value = 2.45
$ws.Cells.Item(242, 16).Value = 0
$ws.Cells.Item(242, 17).Value = 2.025
$ws.Cells.Item(242, 18).Value = 1.825
$ws.Cells.Item(242, 19).Value = 2.75
$ws.Cells.Item(242, 20).Value = 1.925
$ws.Cells.Item(242, 21).Value = 1.925
$ws.Cells.Item(242, 22).Value = 1.7
$ws.Cells.Item(242, 23).Value = -1
$ws.Cells.Item(242, 24).Value = -1
$ws.Cells.Item(242, 25).Value = 1.025
$ws.Cells.Item(242, 26).Value = -1
$ws.Cells.Item(242, 27).Value = -1
$ws.Cells.Item(242, 28).Value = 0.925

$ws.Cells.Item(244, 2).Value = 6870199
$ws.Cells.Item(244, 5).Value = "Cercle Brugge"
$ws.Cells.Item(244, 6).Value = "RWD Molenbeek"
$ws.Cells.Item(244, 7).Value = 4
$ws.Cells.Item(244, 8).Value = 0
$ws.Cells.Item(244, 9).Value = "H"
$ws.Cells.Item(244, 10).Value = 1.363
$ws.Cells.Item(244, 11).Value = 5.5
$ws.Cells.Item(244, 12).Value = 7.5
$ws.Cells.Item(244, 13).Value = 1.3
$ws.Cells.Item(244, 14).Value = 6
$ws.Cells.Item(244, 15).Value = 8.5
$ws.Cells.Item(244, 16).Value = -1.75
$ws.Cells.Item(244, 17).Value = 2.025
$ws.Cells.Item(244, 18).Value = 1.825
$ws.Cells.Item(244, 19).Value = 3.25
$ws.Cells.Item(244, 20).Value = 1.95
$ws.Cells.Item(244, 21).Value = 1.9
$ws.Cells.Item(244, 22).Value = 0.3
$ws.Cells.Item(244, 23).Value = -1
$ws.Cells.Item(244, 24).Value = -1
$ws.Cells.Item(244, 25).Value = 1.025
$ws.Cells.Item(244, 26).Value = -1
$ws.Cells.Item(244, 27).Value = 0.95
$ws.Cells.Item(244, 28).Value = -1

$ws.Cells.Item(275, 2).Value = 7979470
$ws.Cells.Item(275, 5).Value = "Westerlo"
$ws.Cells.Item(275, 6).Value = "OH Leuven"
$ws.Cells.Item(275, 7).Value = 1
$ws.Cells.Item(275, 8).Value = 1
$ws.Cells.Item(275, 9).Value = "D"
$ws.Cells.Item(275, 10).Value = 2.5
$ws.Cells.Item(275, 11).Value = 3.6
$ws.Cells.Item(275, 12).Value = 2.6
$ws.Cells.Item(275, 13).Value = 2.45
$ws.Cells.Item(275, 14).Value = 3.75
$ws.Cells.Item(275, 15).Value = 2.55
$ws.Cells.Item(275, 16).Value = 0
$ws.Cells.Item(275, 17).Value = 1.875
$ws.Cells.Item(275, 18).Value = 1.975
$ws.Cells.Item(275, 19).Value = 3
$ws.Cells.Item(275, 20).Value = 1.85
$ws.Cells.Item(275, 21).Value = 2
$ws.Cells.Item(275, 22).Value = -1
$ws.Cells.Item(275, 23).Value = 2.75
$ws.Cells.Item(275, 24).Value = -1
$ws.Cells.Item(275, 25).Value = 0
$ws.Cells.Item(275, 26).Value = 0
$ws.Cells.Item(275, 27).Value = -1
$ws.Cells.Item(275, 28).Value = 1

$ws.Cells.Item(277, 2).Value = 7979346
$ws.Cells.Item(277, 5).Value = "SintTruidense"
$ws.Cells.Item(277, 6).Value = "Gent"
$ws.Cells.Item(277, 7).Value = 0
$ws.Cells.Item(277, 8).Value = 2
$ws.Cells.Item(277, 9).Value = "A"
$ws.Cells.Item(277, 10).Value = 3.6
$ws.Cells.Item(277, 11).Value = 3.5
$ws.Cells.Item(277, 12).Value = 2
$ws.Cells.Item(277, 13).Value = 3.3
$ws.Cells.Item(277, 14).Value = 3.6
$ws.Cells.Item(277, 15).Value = 2.05
$ws.Cells.Item(277, 16).Value = 0.25
$ws.Cells.Item(277, 17).Value = 2.025
$ws.Cells.Item(277, 18).Value = 1.825
$ws.Cells.Item(277, 19).Value = 3
$ws.Cells.Item(277, 20).Value = 1.975
$ws.Cells.Item(277, 21).Value = 1.875
$ws.Cells.Item(277, 22).Value = -1
$ws.Cells.Item(277, 23).Value = -1
$ws.Cells.Item(277, 24).Value = 1.05
$ws.Cells.Item(277, 25).Value = -1
$ws.Cells.Item(277, 26).Value = 0.825
$ws.Cells.Item(277, 27).Value = -1
$ws.Cells.Item(277, 28).Value = 0.875

$ws.Cells.Item(278, 2).Value = 7979473
$ws.Cells.Item(278, 5).Value = "Anderlecht"
$ws.Cells.Item(278, 6).Value = "Cercle Brugge"
$ws.Cells.Item(278, 7).Value = 3
$ws.Cells.Item(278, 8).Value = 0
$ws.Cells.Item(278, 9).Value = "H"
$ws.Cells.Item(278, 10).Value = 1.909
$ws.Cells.Item(278, 11).Value = 3.6
$ws.Cells.Item(278, 12).Value = 3.8
$ws.Cells.Item(278, 13).Value = 1.8
$ws.Cells.Item(278, 14).Value = 3.8
$ws.Cells.Item(278, 15).Value = 4
$ws.Cells.Item(278, 16).Value = -0.5
$ws.Cells.Item(278, 17).Value = 1.85
$ws.Cells.Item(278, 18).Value = 2
$ws.Cells.Item(278, 19).Value = 2.75
$ws.Cells.Item(278, 20).Value = 1.85
$ws.Cells.Item(278, 21).Value = 2
$ws.Cells.Item(278, 22).Value = 0.8
$ws.Cells.Item(278, 23).Value = -1
$ws.Cells.Item(278, 24).Value = -1
$ws.Cells.Item(278, 25).Value = 0.8500000000000001
$ws.Cells.Item(278, 26).Value = -1
$ws.Cells.Item(278, 27).Value = 0.425
$ws.Cells.Item(278, 28).Value = -0.5

$ws.Cells.Item(279, 2).Value = 7979357
$ws.Cells.Item(279, 5).Value = "Club Brugge"
$ws.Cells.Item(279, 6).Value = "Genk"
$ws.Cells.Item(279, 7).Value = 4
$ws.Cells.Item(279, 8).Value = 0
$ws.Cells.Item(279, 9).Value = "H"
$ws.Cells.Item(279, 10).Value = 1.85
$ws.Cells.Item(279, 11).Value = 3.75
$ws.Cells.Item(279, 12).Value = 3.9
$ws.Cells.Item(279, 13).Value = 1.75
$ws.Cells.Item(279, 14).Value = 3.75
$ws.Cells.Item(279, 15).Value = 4.5
$ws.Cells.Item(279, 16).Value = -0.75
$ws.Cells.Item(279, 17).Value = 2
$ws.Cells.Item(279, 18).Value = 1.85
$ws.Cells.Item(279, 19).Value = 2.75
$ws.Cells.Item(279, 20).Value = 2.025
$ws.Cells.Item(279, 21).Value = 1.825
$ws.Cells.Item(279, 22).Value = 0.75
$ws.Cells.Item(279, 23).Value = -1
$ws.Cells.Item(279, 24).Value = -1
$ws.Cells.Item(279, 25).Value = 1
$ws.Cells.Item(279, 26).Value = -1
$ws.Cells.Item(279, 27).Value = 1.025
$ws.Cells.Item(279, 28).Value = -1

$ws.Cells.Item(297, 17).Value = 1.975
$ws.Cells.Item(297, 18).Value = 1.875
$ws.Cells.Item(297, 20).Value = 1.825
$ws.Cells.Item(297, 21).Value = 2.025

$ws.Cells.Item(299, 13).Value = 2.25
$ws.Cells.Item(299, 14).Value = 3.4
$ws.Cells.Item(299, 17).Value = 1.975
$ws.Cells.Item(299, 18).Value = 1.875

$ws.Cells.Item(300, 13).Value = 1.95
$ws.Cells.Item(300, 15).Value = 3.6
$ws.Cells.Item(300, 17).Value = 2
$ws.Cells.Item(300, 18).Value = 1.85
$ws.Cells.Item(300, 20).Value = 1.925
$ws.Cells.Item(300, 21).Value = 1.925

$ws.Cells.Item(302, 20).Value = 1.975
$ws.Cells.Item(302, 21).Value = 1.875

$ws.Cells.Item(303, 20).Value = 1.95
$ws.Cells.Item(303, 21).Value = 1.9
